$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myoc"
$ws.Range("C2").Value = "Fzd10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.11167
$ws.Range("H2").Value = 3.33501
$ws.Range("I2").Value = 0.1644205457720283
$ws.Range("J2").Value = 0.1644205457720283
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005445666666666667
$ws.Range("N2").Value = 0.016337
$ws.Range("O2").Value = 0.1242225162340131
$ws.Range("P2").Value = 0.1242225162340131
$ws.Range("Q2").Value = 0.006053784263333334
$ws.Range("R2").Value = 0.05448405837
$ws.Range("S2").Value = 0.02042473391637109
$ws.Range("T2").Value = 0.02042473391637109

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myoc"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.11167
$ws.Range("H3").Value = 3.33501
$ws.Range("I3").Value = 0.1644205457720283
$ws.Range("J3").Value = 0.1644205457720283
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03095266666666667
$ws.Range("N3").Value = 0.092858
$ws.Range("O3").Value = 0.706069315814286
$ws.Range("P3").Value = 0.7060693158142859
$ws.Range("Q3").Value = 0.03440915095333333
$ws.Range("R3").Value = 0.30968235858
$ws.Range("S3").Value = 0.1160923022590675
$ws.Range("T3").Value = 0.1160923022590675

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Myoc"
$ws.Range("C4").Value = "Fzd10"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.11167
$ws.Range("H4").Value = 3.33501
$ws.Range("I4").Value = 0.1644205457720283
$ws.Range("J4").Value = 0.1644205457720283
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.007439666666666666
$ws.Range("N4").Value = 0.022319
$ws.Range("O4").Value = 0.169708167951701
$ws.Range("P4").Value = 0.169708167951701
$ws.Range("Q4").Value = 0.008270454243333333
$ws.Range("R4").Value = 0.07443408819
$ws.Range("S4").Value = 0.02790350959658972
$ws.Range("T4").Value = 0.02790350959658972

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Myoc"
$ws.Range("C5").Value = "Fzd10"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.135293
$ws.Range("H5").Value = 15.405879
$ws.Range("I5").Value = 0.7595308659577722
$ws.Range("J5").Value = 0.7595308659577723
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.005445666666666667
$ws.Range("N5").Value = 0.016337
$ws.Range("O5").Value = 0.1242225162340131
$ws.Range("P5").Value = 0.1242225162340131
$ws.Range("Q5").Value = 0.02796509391366667
$ws.Range("R5").Value = 0.251685845223
$ws.Range("S5").Value = 0.0943508353266734
$ws.Range("T5").Value = 0.09435083532667342

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Myoc"
$ws.Range("C6").Value = "Fzd10"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.135293
$ws.Range("H6").Value = 15.405879
$ws.Range("I6").Value = 0.7595308659577722
$ws.Range("J6").Value = 0.7595308659577723
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03095266666666667
$ws.Range("N6").Value = 0.092858
$ws.Range("O6").Value = 0.706069315814286
$ws.Range("P6").Value = 0.7060693158142859
$ws.Range("Q6").Value = 0.1589510124646667
$ws.Range("R6").Value = 1.430559112182
$ws.Range("S6").Value = 0.5362814388666364
$ws.Range("T6").Value = 0.5362814388666364

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Myoc"
$ws.Range("C7").Value = "Fzd10"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.135293
$ws.Range("H7").Value = 15.405879
$ws.Range("I7").Value = 0.7595308659577722
$ws.Range("J7").Value = 0.7595308659577723
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.007439666666666666
$ws.Range("N7").Value = 0.022319
$ws.Range("O7").Value = 0.169708167951701
$ws.Range("P7").Value = 0.169708167951701
$ws.Range("Q7").Value = 0.03820486815566666
$ws.Range("R7").Value = 0.343843813401
$ws.Range("S7").Value = 0.1288985917644625
$ws.Range("T7").Value = 0.1288985917644625

$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Myoc"
$ws.Range("C8").Value = "Fzd10"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.08681433333333333
$ws.Range("H8").Value = 0.260443
$ws.Range("I8").Value = 0.01284019544244376
$ws.Range("J8").Value = 0.01284019544244376
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.005445666666666667
$ws.Range("N8").Value = 0.016337
$ws.Range("O8").Value = 0.1242225162340131
$ws.Range("P8").Value = 0.1242225162340131
$ws.Range("Q8").Value = 0.0004727619212222222
$ws.Range("R8").Value = 0.004254857291
$ws.Range("S8").Value = 0.001595041386796872
$ws.Range("T8").Value = 0.001595041386796872

$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Myoc"
$ws.Range("C9").Value = "Fzd10"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.08681433333333333
$ws.Range("H9").Value = 0.260443
$ws.Range("I9").Value = 0.01284019544244376
$ws.Range("J9").Value = 0.01284019544244376
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.03095266666666667
$ws.Range("N9").Value = 0.092858
$ws.Range("O9").Value = 0.706069315814286
$ws.Range("P9").Value = 0.7060693158142859
$ws.Range("Q9").Value = 0.002687135121555555
$ws.Range("R9").Value = 0.024184216094
$ws.Range("S9").Value = 0.00906606801096798
$ws.Range("T9").Value = 0.00906606801096798

$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Myoc"
$ws.Range("C10").Value = "Fzd10"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.08681433333333333
$ws.Range("H10").Value = 0.260443
$ws.Range("I10").Value = 0.01284019544244376
$ws.Range("J10").Value = 0.01284019544244376
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.007439666666666666
$ws.Range("N10").Value = 0.022319
$ws.Range("O10").Value = 0.169708167951701
$ws.Range("P10").Value = 0.169708167951701
$ws.Range("Q10").Value = 0.0006458697018888888
$ws.Range("R10").Value = 0.005812827317
$ws.Range("S10").Value = 0.002179086044678911
$ws.Range("T10").Value = 0.002179086044678911

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Myoc"
$ws.Range("C11").Value = "Fzd10"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.34273
$ws.Range("H11").Value = 1.02819
$ws.Range("I11").Value = 0.05069117062837646
$ws.Range("J11").Value = 0.05069117062837647
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.005445666666666667
$ws.Range("N11").Value = 0.016337
$ws.Range("O11").Value = 0.1242225162340131
$ws.Range("P11").Value = 0.1242225162340131
$ws.Range("Q11").Value = 0.001866393336666667
$ws.Range("R11").Value = 0.01679754003
$ws.Range("S11").Value = 0.006296984766304625
$ws.Range("T11").Value = 0.006296984766304625

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Myoc"
$ws.Range("C12").Value = "Fzd10"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.34273
$ws.Range("H12").Value = 1.02819
$ws.Range("I12").Value = 0.05069117062837646
$ws.Range("J12").Value = 0.05069117062837647
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.03095266666666667
$ws.Range("N12").Value = 0.092858
$ws.Range("O12").Value = 0.706069315814286
$ws.Range("P12").Value = 0.7060693158142859
$ws.Range("Q12").Value = 0.01060840744666667
$ws.Range("R12").Value = 0.09547566701999999
$ws.Range("S12").Value = 0.035791480163403
$ws.Range("T12").Value = 0.035791480163403

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Myoc"
$ws.Range("C13").Value = "Fzd10"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.34273
$ws.Range("H13").Value = 1.02819
$ws.Range("I13").Value = 0.05069117062837646
$ws.Range("J13").Value = 0.05069117062837647
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.007439666666666666
$ws.Range("N13").Value = 0.022319
$ws.Range("O13").Value = 0.169708167951701
$ws.Range("P13").Value = 0.169708167951701
$ws.Range("Q13").Value = 0.002549796956666666
$ws.Range("R13").Value = 0.02294817261
$ws.Range("S13").Value = 0.008602705698668844
$ws.Range("T13").Value = 0.008602705698668845

$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Myoc"
$ws.Range("C14").Value = "Fzd10"
$ws.Range("D14").Value = "FAPs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.08463066666666667
$ws.Range("H14").Value = 0.253892
$ws.Range("I14").Value = 0.01251722219937926
$ws.Range("J14").Value = 0.01251722219937926
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.005445666666666667
$ws.Range("N14").Value = 0.016337
$ws.Range("O14").Value = 0.1242225162340131
$ws.Range("P14").Value = 0.1242225162340131
$ws.Range("Q14").Value = 0.0004608704004444445
$ws.Range("R14").Value = 0.004147833604
$ws.Range("S14").Value = 0.001554920837867139
$ws.Range("T14").Value = 0.001554920837867139

$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Myoc"
$ws.Range("C15").Value = "Fzd10"
$ws.Range("D15").Value = "Inflammatory-Mac"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.08463066666666667
$ws.Range("H15").Value = 0.253892
$ws.Range("I15").Value = 0.01251722219937926
$ws.Range("J15").Value = 0.01251722219937926
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.03095266666666667
$ws.Range("N15").Value = 0.092858
$ws.Range("O15").Value = 0.706069315814286
$ws.Range("P15").Value = 0.7060693158142859
$ws.Range("Q15").Value = 0.002619544815111111
$ws.Range("R15").Value = 0.023575903336
$ws.Range("S15").Value = 0.008838026514211105
$ws.Range("T15").Value = 0.008838026514211105

$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Myoc"
$ws.Range("C16").Value = "Fzd10"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.08463066666666667
$ws.Range("H16").Value = 0.253892
$ws.Range("I16").Value = 0.01251722219937926
$ws.Range("J16").Value = 0.01251722219937926
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.007439666666666666
$ws.Range("N16").Value = 0.022319
$ws.Range("O16").Value = 0.169708167951701
$ws.Range("P16").Value = 0.169708167951701
$ws.Range("Q16").Value = 0.0006296239497777779
$ws.Range("R16").Value = 0.005666615547999999
$ws.Range("S16").Value = 0.002124274847301015
$ws.Range("T16").Value = 0.002124274847301015

